# Updated notebook, reran simulation
# - Added two new strings to the shared-string pool ("Holden", "Rizzie Spiral")
# - Renamed "Thomas Hex" -> "Matthies Hex"
# - Inserted two new data rows (new rows 4 & 5) into the simulation results,
#   pushing the previous rows 4-29 down to rows 6-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two blank rows at row 4 (existing rows 4:29 shift down to 6:31)
# ---------------------------------------------------------------------------
$ws.Rows("4:5").Insert(-4142)   # xlFormatFromRightOrBelow - avoid bad style copy

# Column A on the new rows needs the same bold/centered/bordered style used by
# the rest of column A (style gets lost/garbled by the row insert above).
foreach ($r in 4..5) {
    $c = $ws.Cells.Item($r, 1)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 2. Rename the "Thomas Hex" label to "Matthies Hex" (wherever it appears)
# ---------------------------------------------------------------------------
$used = $ws.UsedRange
$found = $used.Find("Thomas Hex")
if ($found) {
    $found.Value2 = "Matthies Hex"
}

# ---------------------------------------------------------------------------
# 3. Fill in the two new data rows (row 4 and row 5)
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value2 = 2
$ws.Cells.Item(4, 2).Value2 = "Holden"

$ws.Cells.Item(5, 1).Value2 = 3
$ws.Cells.Item(5, 2).Value2 = "Rizzie Spiral"

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

$row4vals = @(0.8337253803315133, 0.9558137102140735, 0.8784717762277323, 2.471020677928275, 0.8337253803315133, 0.8533298139341446, 1.115667275503263, 0.878471776227579, 0.878471776227579, 0.8414321326223833, 0.8709370640830574, 0.878471776227579, 2.471020677928275, 1.652373029129894, 1.66217524593121, 1.394405944829122, 1.386025290731311, 1.394405944829122, 1.259136912105378, 1.183003884929818, 1.102549728855536)

$row5vals = @(-0.001352288825052127, 0.01029846540936226, 0.001260769884377864, 4.49693082844336, -0.001352288825052127, 0.1797443978197705, 1.096825101568488, 0.001260769884377864, 0.001260769884377864, -0.0003898151431831477, 2.886504778798488, 0.001260769884377864, 4.49693082844336, 2.247789269809154, 2.338337613131565, 1.498946436500895, 1.558440979146026, 1.498946436500895, 1.169145926830614, 0.9355688954413669, 1.083727779744451)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value2 = $row4vals[$i]
    $ws.Range($cols[$i] + "5").Value2 = $row5vals[$i]
}

# ---------------------------------------------------------------------------
# Sanity: dimension should now read A1:W31 (handled automatically by Excel)
# ---------------------------------------------------------------------------
Write-Host ("New used range: " + $ws.UsedRange.Address())
